$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7500
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 7500
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 22500
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -24248
$ws.Range("H70").Value = 1270.7142
$ws.Range("I70").Value = 697.5
$ws.Range("K70").Value = 2092.5
$ws.Range("M70").Value = -1822.5
$ws.Range("H72").Value = 7500
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 7500
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 67500
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -76236
$ws.Range("H73").Value = 1270.7142
$ws.Range("I73").Value = 697.5
$ws.Range("K73").Value = 2092.5
$ws.Range("M73").Value = -1156.5
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()
$ws.Range("H82").Value = 690
$ws.Range("I82").Value = 690
$ws.Range("K82").Value = 2070
$ws.Range("M82").Value = -1664
$ws.Range("H85").Value = 690
$ws.Range("I85").Value = 690
$ws.Range("K85").Value = 2070
$ws.Range("M85").Value = -666
$ws.Range("H137").Value = 2326.8572
$ws.Range("J137").Value = 2903
$ws.Range("L137").Value = 8709
$ws.Range("N137").Value = -13809
$ws.Range("H138").Value = 2648.6924
$ws.Range("J138").Value = 4299.7144
$ws.Range("L138").Value = 12899.1432
$ws.Range("N138").Value = -23179.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1161.25
$ws.Range("I88").Value = 903
$ws.Range("J88").Value = 1419.5
$ws.Range("K88").Value = 903
$ws.Range("L88").Value = 1419.5
$ws.Range("M88").Value = -497
$ws.Range("N88").Value = -2231.5
$ws.Range("H91").Value = 1161.25
$ws.Range("I91").Value = 903
$ws.Range("J91").Value = 1419.5
$ws.Range("K91").Value = 903
$ws.Range("L91").Value = 1419.5
$ws.Range("M91").Value = 501
$ws.Range("N91").Value = -4227.5
$ws.Range("H108").Value = 31810.5
$ws.Range("I108").Value = 33621
$ws.Range("J108").Value = 30000
$ws.Range("K108").Value = 33621
$ws.Range("L108").Value = 30000
$ws.Range("M108").Value = -29781
$ws.Range("N108").Value = -37680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5767.75
$ws.Range("I86").Value = 4424
$ws.Range("K86").Value = 4424
$ws.Range("M86").Value = -3301
$ws.Range("H89").Value = 5767.75
$ws.Range("I89").Value = 4424
$ws.Range("K89").Value = 22120
$ws.Range("M89").Value = -16504
$ws.Range("H107").Value = 1288.3334
$ws.Range("J107").Value = 1288.3334
$ws.Range("L107").Value = 1288.3334
$ws.Range("N107").Value = -5128.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 27000
$ws.Range("J68").Value = 27000
$ws.Range("L68").Value = 27000
$ws.Range("N68").Value = -28498
$ws.Range("H71").Value = 27000
$ws.Range("J71").Value = 27000
$ws.Range("L71").Value = 81000
$ws.Range("N71").Value = -88488
$ws.Range("H99").Value = 7000
$ws.Range("J99").Value = 6500
$ws.Range("L99").Value = 6500
$ws.Range("N99").Value = -9496
$ws.Range("H122").Value = 1824
$ws.Range("I122").Value = 1824
$ws.Range("K122").Value = 5472
$ws.Range("M122").Value = -3022
$ws.Range("H126").Value = 7000
$ws.Range("J126").Value = 6500
$ws.Range("L126").Value = 19500
$ws.Range("N126").Value = -24440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 338.33334
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 338.33334
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 1015.00002
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -1353.00002
$ws.Range("H30").Value = 338.33334
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 338.33334
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 1015.00002
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -1219.00002
$ws.Range("H54").Value = 1655.8
$ws.Range("I54").Value = 1004
$ws.Range("J54").Value = 1818.75
$ws.Range("K54").Value = 3012
$ws.Range("L54").Value = 5456.25
$ws.Range("M54").Value = -2453
$ws.Range("N54").Value = -6574.25
$ws.Range("H55").Value = 2764
$ws.Range("J55").Value = 3750
$ws.Range("L55").Value = 11250
$ws.Range("N55").Value = -11604
$ws.Range("H134").Value = 2515
$ws.Range("I134").Value = 30
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 90
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = 4980
$ws.Range("N134").Value = -25140
$ws.Range("H138").Value = 2000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 6000
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -16280
$ws.Range("H139").Value = 1819.4286
$ws.Range("I139").Value = 1456
$ws.Range("K139").Value = 4368
$ws.Range("M139").Value = 772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4475
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H130").Value = 120000
$ws.Range("J130").Value = 120000
$ws.Range("L130").Value = 120000
$ws.Range("N130").Value = -130040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
